$wb = $excel.ActiveWorkbook

# The same F3:F9 "想去人数" (attendance count) update applies to both the
# "展览" sheet and the "全部类型" sheet, which mirror each other.
$sheetNames = @("展览", "全部类型")

# Mapping of row number -> new value for column F
$updates = @{
    3 = 1823
    4 = 352
    5 = 1123
    6 = 1040
    7 = 45
    8 = 5919
    9 = 95
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
